$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.774.66'
$ws.Range("E2").Value = '  +7.33%  '
$ws.Range("D3").Value = '1.951.02'
$ws.Range("E3").Value = '  +5.55%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  -0.40%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '341.78'
$ws.Range("E5").Value = '  +2.04%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  -0.31%  '
$ws.Range("E7").Value = '  +2.76%  '
$ws.Range("E8").Value = '  +7.31%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '47.86'
$ws.Range("E9").Value = '  +2.56%  '
$ws.Range("E10").Value = '  +4.32%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.038'
$ws.Range("E11").Value = '  +7.00%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '22.76'
$ws.Range("E12").Value = '  +6.92%  '
$ws.Range("D13").Value = '1.950.38'
$ws.Range("E13").Value = '  +6.72%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.159'
$ws.Range("E14").Value = '  +4.58%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.389'
$ws.Range("E15").Value = '  +3.46%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '91.80'
$ws.Range("E16").Value = '  +2.10%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.002'
$ws.Range("E17").Value = '  -0.36%  '
$ws.Range("E18").Value = '  +3.03%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06688'
$ws.Range("E19").Value = '  +1.17%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.06'
$ws.Range("E20").Value = '  +4.15%  '
$ws.Range("E21").Value = '  -0.42%  '
$ws.Range("D22").Value = '29.737.98'
$ws.Range("E22").Value = '  +7.36%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.600'
$ws.Range("E23").Value = '  +4.54%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.27'
$ws.Range("E24").Value = '  +3.89%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.276'
$ws.Range("E25").Value = '  -0.39%  '
$ws.Range("E26").Value = '  +6.59%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '160.97'
$ws.Range("E27").Value = '  +1.29%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.26'
$ws.Range("E28").Value = '  +3.90%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.173'
$ws.Range("E29").Value = '  +4.91%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.677'
$ws.Range("E30").Value = '  +6.22%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '122.91'
$ws.Range("E31").Value = '  +3.60%  '
$ws.Range("E32").Value = '  +6.34%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09629'
$ws.Range("E33").Value = '  +2.23%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.476'
$ws.Range("E34").Value = '  +11.02%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.681'
$ws.Range("E35").Value = '  +2.68%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.507'
$ws.Range("E36").Value = '  +4.65%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06267'
$ws.Range("E37").Value = '  +4.49%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02313'
$ws.Range("E38").Value = '  +4.65%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.483'
$ws.Range("E39").Value = '  +3.08%  '
$ws.Range("B40").Value = 'TrustWalletToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.186'
$ws.Range("E40").Value = '  +2.75%  '
$ws.Range("B41").Value = 'TheSandbox'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6097'
$ws.Range("E41").Value = '  +5.19%  '
$ws.Range("E42").Value = '  +6.95%  '
$ws.Range("E43").Value = '  -0.27%  '
$ws.Range("E44").Value = '  +2.99%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.277'
$ws.Range("E45").Value = '  -0.14%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.389'
$ws.Range("E46").Value = '  +32.34%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5715'
$ws.Range("E47").Value = '  +5.02%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '12.50'
$ws.Range("E48").Value = '  +5.17%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.07415'
$ws.Range("E49").Value = '  +8.34%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.991'
$ws.Range("E50").Value = '  +3.33%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '112.94'
$ws.Range("E51").Value = '  +1.78%  '
